$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a two-row header (row 1 + row 2) describing the
# plant table, followed by two data rows (row 3: Kubel, row 4: Wasserauen)
# and a block of empty filler rows down to row 33.
#
# The target layout collapses the header into a single row 1 with new
# column titles (idx, idx2, Name, Date Start, Date End, (m3/s), (MW1),
# (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year), and the two data rows
# shift up to rows 2 and 3. That is exactly what deleting the current
# row 2 (the secondary header row) accomplishes, since it shifts rows
# 3-33 up to 2-32, and the trailing filler row count naturally shrinks
# from 29 (rows 5-33) to 29 (rows 4-32) - i.e. one fewer total row.
$ws.Rows.Item(2).Delete()

# Create a temporary named cell style that mirrors the existing "(m3/s)"
# style (Arial 9) but - unlike the current style 2 - does NOT mark the
# (default, General) number format as explicitly applied. Applying it to
# a cell bakes a matching cell format (xf) into the workbook's styles,
# after which we can discard the temporary named style again, leaving
# only the new low-level cell format behind (re-pointed at the default
# cellStyleXfs entry), matching the new cellXfs entry added upstream.
$tempStyleName = "TempHeaderStyle"
$tempStyle = $wb.Styles.Add($tempStyleName)
$tempStyle.Font.Name = "Arial"
$tempStyle.Font.Size = 9
$tempStyle.IncludeNumber = $false

# Row 1: brand new header row. Columns A-E get the plain default style;
# E1 previously held a styled header ("mation" truncation, style 2) so
# its format must be reset back to the workbook default explicitly.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("E1").Style = "Normal"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("F1").Style = $tempStyleName

$ws.Range("G1").Value = "(MW1)"
$ws.Range("G1").Style = $tempStyleName

$ws.Range("H1").Value = "(MW2)"
$ws.Range("H1").Style = $tempStyleName

$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("I1").Style = $tempStyleName

$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("J1").Style = $tempStyleName

$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("K1").Style = $tempStyleName

# Remove the temporary named style now that its format has been baked
# into the cells that used it; this keeps cellStyles/cellStyleXfs at
# their original count of 3.
$wb.Styles($tempStyleName).Delete()

# Match the selection left behind by the edit (the data rows, A2:K2).
$ws.Range("A2:K2").Select()
